$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("L1").Value = "status_validacao"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Status values for rows 2-21
$statuses = @("OK", "OK", "OK", "OK", "OK", "ERRO", "OK", "OK", "OK", "OK", "ERRO", "ERRO", "ERRO", "ERRO", "ERRO", "ERRO", "ERRO", "ERRO", "ERRO", "ERRO")

for ($i = 0; $i -lt $statuses.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $statuses[$i]
}
